$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell F1, styled like the other headers (copy style from E1)
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$timestamps = @(
    "2021-10-05 13:40:05.037318",
    "2021-10-05 13:40:05.037334",
    "2021-10-05 13:40:05.037338",
    "2021-10-05 13:40:05.037341",
    "2021-10-05 13:40:05.037345",
    "2021-10-05 13:40:05.037348",
    "2021-10-05 13:40:05.037366",
    "2021-10-05 13:40:05.037368",
    "2021-10-05 13:40:05.037371",
    "2021-10-05 13:40:05.037374",
    "2021-10-05 13:40:05.037376",
    "2021-10-05 13:40:05.037379"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
